$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column G
$wsOverview.Range("G2").Value = "2016-09-05 12:20:20"
$wsOverview.Range("G5").Value = "2016-09-05 12:20:20"

# zh-cn sheet: Status column E -> "mt"
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime column H
$wsZhCn.Range("H2").Value = "2016-09-05 12:20:01"
$wsZhCn.Range("H5").Value = "2016-09-05 12:20:01"

# zh-cn sheet: Correspond Handback DateTime column K
$wsZhCn.Range("K2").Value = "2016-09-05 12:20:37"
$wsZhCn.Range("K5").Value = "2016-09-05 12:20:37"

# de-de sheet: Status column E -> "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# de-de sheet: Correspond Handoff Datetime column H (shares value w/ Overview's date)
$wsDeDe.Range("H2").Value = "2016-09-05 12:20:20"
$wsDeDe.Range("H5").Value = "2016-09-05 12:20:20"

# de-de sheet: Correspond Handback DateTime column K
$wsDeDe.Range("K2").Value = "2016-09-05 12:20:44"
$wsDeDe.Range("K5").Value = "2016-09-05 12:20:44"
